$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Populate new cells in the order the shared-string table records them so
# that the resulting xl/sharedStrings.xml index order matches the target:
#   6  Sacaremos el promedio de las edades de los alumnos  (C11)
#   7  lugar de residencia                                  (E4)
#   8  lima                                                  (E5 / E6 / H15)
#   9  cañete                                                (E7)
#  10  ica                                                   (E8)
#  11  la moda de lugar de residencia es  :                  (C15)
# ---------------------------------------------------------------------------

$ws.Range("C11").Value = "Sacaremos el promedio de las edades de los alumnos"

$ws.Range("E4").Value = "lugar de residencia"
$ws.Range("E5").Value = "lima"
$ws.Range("E6").Value = "lima"
$ws.Range("E7").Value = "cañete"
$ws.Range("E8").Value = "ica"

$ws.Range("C15").Value = "la moda de lugar de residencia es  :"

# Formula: average of the four ages, placed in H11
$ws.Range("H11").Formula = "=SUM(D5+D6+D7+D8)/4"

# H15 holds the mode of "lugar de residencia" (reuses the "lima" shared string)
$ws.Range("H15").Value = "lima"

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# E4:E8 get the same thin-border style already used by C5:D8 (cellXfs index 1).
# Copying the format from an existing bordered cell reuses that style instead
# of fabricating a new (slightly different) border definition.
$ws.Range("C5").Copy()
$ws.Range("E4:E8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# H11 and H15 get a solid yellow fill (new cellXfs entry, fillId -> FFFFFF00).
$ws.Range("H11").Interior.Color = 65535
$ws.Range("H15").Interior.Color = 65535

# Column widths for the newly used columns E and H.
$ws.Columns.Item(5).ColumnWidth = 17.6
$ws.Columns.Item(8).ColumnWidth = 11

# Selection moves to K7, as recorded in the saved sheet view.
$ws.Range("K7").Select()
